$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: rows 8-13 ("Ready for handoff" files) are being handed off now.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E8:E13").Value = "ht"
$wsZh.Range("H8:H13").Value = "2016-09-02 20:25:30"

# "de-de" sheet: same rows/files, different handoff timestamp.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E8:E13").Value = "ht"
$wsDe.Range("H8:H13").Value = "2016-09-02 20:25:36"

# "Overview" sheet mirrors the de-de handoff datetime in column G for the
# same rows; keep it in sync with the de-de sheet's new timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G8:G13").Value = "2016-09-02 20:25:36"
